$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.544.36"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "2.603.37"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.07"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.79"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "3.061.38"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "59.450.57"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.81"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.586.20"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.43"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.10"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.23"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.68"
$ws.Range("E30").Value = "  +6.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.89"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.92"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.849"
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.827"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "273.23"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.72"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.57"
$ws.Range("E46").Value = "  +3.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0224"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "1.940.91"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.05"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("E51").Value = "  +1.61%  "
